# Auto-generated Excel COM-interop script
# Updates market-price snapshot values (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, matching the scheduled-runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3939.6667
$ws.Range("I62").Value = 3586.3333
$ws.Range("K62").Value = 3586.3333
$ws.Range("M62").Value = -2962.3333
$ws.Range("H64").Value = 4572.2
$ws.Range("I64").Value = 4487.5
$ws.Range("J64").Value = 4628.6665
$ws.Range("K64").Value = 4487.5
$ws.Range("L64").Value = 4628.6665
$ws.Range("M64").Value = -4239.5
$ws.Range("N64").Value = -5124.6665
$ws.Range("H65").Value = 3939.6667
$ws.Range("I65").Value = 3586.3333
$ws.Range("K65").Value = 17931.6665
$ws.Range("M65").Value = -14811.6665
$ws.Range("H67").Value = 4572.2
$ws.Range("I67").Value = 4487.5
$ws.Range("J67").Value = 4628.6665
$ws.Range("K67").Value = 4487.5
$ws.Range("L67").Value = 4628.6665
$ws.Range("M67").Value = -3629.5
$ws.Range("N67").Value = -6344.6665
$ws.Range("H70").Value = 3636.4707
$ws.Range("I70").Value = 2337.4
$ws.Range("K70").Value = 7012.200000000001
$ws.Range("M70").Value = -6742.200000000001
$ws.Range("H73").Value = 3636.4707
$ws.Range("I73").Value = 2337.4
$ws.Range("K73").Value = 7012.200000000001
$ws.Range("M73").Value = -6076.200000000001
$ws.Range("H74").Value = 48708.453
$ws.Range("I74").Value = 3486.625
$ws.Range("J74").Value = 169300
$ws.Range("K74").Value = 3486.625
$ws.Range("L74").Value = 169300
$ws.Range("M74").Value = -2550.625
$ws.Range("N74").Value = -171172
$ws.Range("H77").Value = 48708.453
$ws.Range("I77").Value = 3486.625
$ws.Range("J77").Value = 169300
$ws.Range("K77").Value = 17433.125
$ws.Range("L77").Value = 846500
$ws.Range("M77").Value = -12753.125
$ws.Range("N77").Value = -855860
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H100").Value = 2956.5715
$ws.Range("I100").Value = 2959.2
$ws.Range("K100").Value = 2959.2
$ws.Range("M100").Value = -2418.2
$ws.Range("H103").Value = 1144.7778
$ws.Range("J103").Value = 1190.6
$ws.Range("L103").Value = 3571.8
$ws.Range("N103").Value = -4743.799999999999
$ws.Range("H111").Value = 500
$ws.Range("I111").Value = 450
$ws.Range("K111").Value = 1350
$ws.Range("M111").Value = 1717
$ws.Range("H132").Value = 2052.9
$ws.Range("I132").Value = 2052.9
$ws.Range("K132").Value = 6158.700000000001
$ws.Range("M132").Value = -3628.700000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1158.8
$ws.Range("J2").Value = 1099.5
$ws.Range("L2").Value = 1099.5
$ws.Range("N2").Value = -1325.5
$ws.Range("H45").Value = 1282.5
$ws.Range("I45").Value = 1282.5
$ws.Range("K45").Value = 1282.5
$ws.Range("M45").Value = -905.5
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H88").Value = 1978.4706
$ws.Range("J88").Value = 2096.3333
$ws.Range("L88").Value = 2096.3333
$ws.Range("N88").Value = -2908.3333
$ws.Range("H91").Value = 1978.4706
$ws.Range("J91").Value = 2096.3333
$ws.Range("L91").Value = 2096.3333
$ws.Range("N91").Value = -4904.3333
$ws.Range("H116").Value = 1158.8
$ws.Range("J116").Value = 1099.5
$ws.Range("L116").Value = 1099.5
$ws.Range("N116").Value = -5687.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1158.8
$ws.Range("J3").Value = 1099.5
$ws.Range("L3").Value = 1099.5
$ws.Range("N3").Value = -1327.5
$ws.Range("H11").Value = 345.5
$ws.Range("I11").Value = 345.5
$ws.Range("K11").Value = 345.5
$ws.Range("M11").Value = -205.5
$ws.Range("H20").Value = 4474
$ws.Range("I20").Value = 3969
$ws.Range("K20").Value = 3969
$ws.Range("M20").Value = -3722
$ws.Range("H94").Value = 1402
$ws.Range("H105").Value = 3819.75
$ws.Range("I105").Value = 3651.2856
$ws.Range("K105").Value = 3651.2856
$ws.Range("M105").Value = -1904.2856
$ws.Range("H134").Value = 13006
$ws.Range("I134").Value = 13008.333
$ws.Range("K134").Value = 39024.999
$ws.Range("M134").Value = -36489.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3223.4614
$ws.Range("I58").Value = 3056.111
$ws.Range("K58").Value = 3056.111
$ws.Range("M58").Value = -2853.111
$ws.Range("H86").Value = 7138.6
$ws.Range("I86").Value = 6736.1113
$ws.Range("K86").Value = 6736.1113
$ws.Range("M86").Value = -5613.1113
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H89").Value = 7138.6
$ws.Range("I89").Value = 6736.1113
$ws.Range("K89").Value = 33680.5565
$ws.Range("M89").Value = -28064.5565
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H103").Value = 5524
$ws.Range("I103").Value = 5524
$ws.Range("K103").Value = 5524
$ws.Range("M103").Value = -4352
$ws.Range("H108").Value = 38396
$ws.Range("H136").Value = 3223.4614
$ws.Range("I136").Value = 3056.111
$ws.Range("K136").Value = 9168.332999999999
$ws.Range("M136").Value = -6618.332999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 316.125
$ws.Range("I2").Value = 397.66666
$ws.Range("J2").Value = 267.2
$ws.Range("K2").Value = 2385.99996
$ws.Range("L2").Value = 1603.2
$ws.Range("M2").Value = -2272.99996
$ws.Range("N2").Value = -1829.2
$ws.Range("H12").Value = 165.75
$ws.Range("J12").Value = 211
$ws.Range("L12").Value = 633
$ws.Range("N12").Value = -979
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3338
$ws.Range("H26").Value = 72914.42999999999
$ws.Range("I26").Value = 143642.72
$ws.Range("J26").Value = 2186.1428
$ws.Range("K26").Value = 430928.16
$ws.Range("L26").Value = 6558.428400000001
$ws.Range("M26").Value = -430640.16
$ws.Range("N26").Value = -7134.428400000001
$ws.Range("H27").Value = 1000
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3204
$ws.Range("H29").Value = 244.6
$ws.Range("I29").Value = 207.66667
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 623.00001
$ws.Range("L29").Value = 900
$ws.Range("M29").Value = -346.00001
$ws.Range("N29").Value = -1454
$ws.Range("H37").Value = 99712.71000000001
$ws.Range("J37").Value = 99712.71000000001
$ws.Range("L37").Value = 299138.13
$ws.Range("N37").Value = -299362.13

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 124.55556
$ws.Range("I2").Value = 99.8
$ws.Range("K2").Value = 99.8
$ws.Range("M2").Value = 13.2
$ws.Range("H70").Value = 1669
$ws.Range("J70").Value = 2115
$ws.Range("L70").Value = 2115
$ws.Range("N70").Value = -2655
$ws.Range("H73").Value = 1669
$ws.Range("J73").Value = 2115
$ws.Range("L73").Value = 2115
$ws.Range("N73").Value = -3987
$ws.Range("H80").Value = 2621.75
$ws.Range("I80").Value = 2245
$ws.Range("J80").Value = 2998.5
$ws.Range("K80").Value = 2245
$ws.Range("L80").Value = 2998.5
$ws.Range("M80").Value = -1247
$ws.Range("N80").Value = -4994.5
$ws.Range("H83").Value = 2621.75
$ws.Range("I83").Value = 2245
$ws.Range("J83").Value = 2998.5
$ws.Range("K83").Value = 11225
$ws.Range("L83").Value = 14992.5
$ws.Range("M83").Value = -6233
$ws.Range("N83").Value = -24976.5
$ws.Range("H100").Value = 38950
$ws.Range("J100").Value = 38950
$ws.Range("L100").Value = 38950
$ws.Range("N100").Value = -41114
$ws.Range("H102").Value = 3299.375
$ws.Range("I102").Value = 1628.1428
$ws.Range("J102").Value = 14998
$ws.Range("K102").Value = 1628.1428
$ws.Range("L102").Value = 14998
$ws.Range("M102").Value = -6.142800000000079
$ws.Range("N102").Value = -18242
$ws.Range("H126").Value = 7439.8
$ws.Range("I126").Value = 11666.333
$ws.Range("K126").Value = 34998.999
$ws.Range("M126").Value = -32528.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2128
$ws.Range("I68").Value = 2256
$ws.Range("K68").Value = 2256
$ws.Range("M68").Value = -1507
$ws.Range("H71").Value = 2128
$ws.Range("I71").Value = 2256
$ws.Range("K71").Value = 11280
$ws.Range("M71").Value = -7536
$ws.Range("H136").Value = 2877.7856
$ws.Range("I136").Value = 2690.3333
$ws.Range("K136").Value = 8070.999899999999
$ws.Range("M136").Value = -5520.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H117").Value = 20000
$ws.Range("J117").Value = 20000
$ws.Range("L117").Value = 20000
$ws.Range("N117").Value = -29178
